{"js": "// Locate the three \"Steps:\" list paragraphs we need to edit by their\n// (pre-edit) text content.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (const p of paragraphs.items) {\n  p.load(\"text\");\n}\nawait context.sync();\n\nfunction findParagraph(items, needle) {\n  for (const p of items) {\n    if (p.text.includes(needle)) {\n      return p;\n    }\n  }\n  return null;\n}\n\nconst p1 = findParagraph(paragraphs.items, \"Get obs flow data from wdm file using tsproc\");\nconst p2 = findParagraph(paragraphs.items, \"Remove flows from obs data used for calibration.\");\nconst p3 = findParagraph(paragraphs.items, \"Process obs flow with data removed for pest run\");\n\nif (!p1) throw new Error(\"Could not locate step 1 paragraph\");\nif (!p2) throw new Error(\"Could not locate step 2 paragraph\");\nif (!p3) throw new Error(\"Could not locate step 3 paragraph\");\n\n// Shared OOXML wrapper for Paragraph.insertOoxml(..., Word.InsertLocation.replace).\n// Rebuilding the paragraph from raw WordprocessingML (rather than just\n// setting .text) lets us reproduce the exact run layout the diff shows:\n// <w:proofErr/> spell-check markers bracketing \"obs\"/\"wdm\"/\"hspf\", and a\n// manual <w:br/> line break before the r-script path.\nfunction pkgXml(bodyInner) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    \"<w:body>\" +\n    bodyInner +\n    \"</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>\"\n  );\n}\n\nconst listPPr =\n  '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>';\n\n// --- Step 1: \"Get obs flow data from wdm file using tsproc\"\n//            -> \"Get obs flow data from wdm file using hspf\" (+ proofErr marks)\nconst body1 =\n  \"<w:p>\" +\n  listPPr +\n  '<w:r><w:t xml:space=\"preserve\">Get </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>obs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> flow data from </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>wdm</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> file using </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>hspf</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  \"</w:p>\";\np1.insertOoxml(pkgXml(body1), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 2: \"Remove flows from obs data used for calibration.\" (text unchanged, proofErr added)\nconst body2 =\n  \"<w:p>\" +\n  listPPr +\n  '<w:r><w:t xml:space=\"preserve\">Remove flows from </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>obs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> data used for calibration.</w:t></w:r>' +\n  \"</w:p>\";\np2.insertOoxml(pkgXml(body2), Word.InsertLocation.replace);\nawait context.sync();\n\n// --- Step 3: \"Process obs flow with data removed for pest run\"\n//            -> append r-script pointer on a new line within the same paragraph.\n// p3 is the last paragraph in the document body, and insertOoxml(...,\n// replace) cannot consume the body's final paragraph mark, so replacing it\n// in place would leave a stray empty paragraph behind. Work around that by\n// temporarily adding a throw-away paragraph after it (so our target is no\n// longer \"last\"), doing the OOXML replace, then deleting the throw-away\n// paragraph, which ends up as the new (empty) final paragraph.\np3.load(\"isLastParagraph\");\nawait context.sync();\nconst wasLastParagraph = p3.isLastParagraph;\nif (wasLastParagraph) {\n  p3.insertParagraph(\"\", Word.InsertLocation.after);\n  await context.sync();\n}\n\nconst body3 =\n  \"<w:p>\" +\n  listPPr +\n  '<w:r><w:t xml:space=\"preserve\">Process </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>obs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t xml:space=\"preserve\"> flow with data removed for pest run</w:t></w:r>' +\n  \"<w:r><w:t>, use previously developed r-script:</w:t></w:r>\" +\n  \"<w:r><w:br/></w:r>\" +\n  \"<w:r><w:t>M:\\\\Models\\\\Bacteria\\\\HSPF\\\\bigelkHydroCal201601\\\\r-files</w:t></w:r>\" +\n  \"<w:r><w:t>\\\\hspf.output-proc.R</w:t></w:r>\" +\n  \"</w:p>\";\np3.insertOoxml(pkgXml(body3), Word.InsertLocation.replace);\nawait context.sync();\n\nif (wasLastParagraph) {\n  const paragraphsAfter = body.paragraphs;\n  paragraphsAfter.load(\"items\");\n  await context.sync();\n  const trailing = paragraphsAfter.items[paragraphsAfter.items.length - 1];\n  trailing.delete();\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Namespaced OOXML package wrapper used with Range.InsertXML to rebuild a\n# single paragraph's run content (this lets us place <w:proofErr/> markers\n# and a <w:br/> exactly where Word's editor would have put them).\nfunction New-PkgXml([string]$bodyInner) {\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# Shared list-paragraph properties (unchanged by the edit).\n$listPPr = '<w:pPr><w:pStyle w:val=\"ListParagraph\"/><w:numPr><w:ilvl w:val=\"0\"/><w:numId w:val=\"1\"/></w:numPr></w:pPr>'\n\nfunction Find-Paragraph([string]$needle) {\n    foreach ($p in $d.Paragraphs) {\n        if ($p.Range.Text.Contains($needle)) {\n            return $p\n        }\n    }\n    return $null\n}\n\n# --- Step 1: \"Get obs flow data from wdm file using tsproc\"\n#            -> \"Get obs flow data from wdm file using hspf\" (+ proofErr marks)\n$p1 = Find-Paragraph \"Get obs flow data from wdm file using tsproc\"\nif ($null -eq $p1) { throw \"Could not locate step 1 paragraph\" }\n$body1 = '<w:p>' + $listPPr +\n    '<w:r><w:t xml:space=\"preserve\">Get </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>obs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> flow data from </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>wdm</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> file using </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>hspf</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '</w:p>'\n$p1.Range.InsertXML((New-PkgXml $body1))\n\n# --- Step 2: \"Remove flows from obs data used for calibration.\" (text unchanged, proofErr added)\n$p2 = Find-Paragraph \"Remove flows from obs data used for calibration.\"\nif ($null -eq $p2) { throw \"Could not locate step 2 paragraph\" }\n$body2 = '<w:p>' + $listPPr +\n    '<w:r><w:t xml:space=\"preserve\">Remove flows from </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>obs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> data used for calibration.</w:t></w:r>' +\n    '</w:p>'\n$p2.Range.InsertXML((New-PkgXml $body2))\n\n# --- Step 3: \"Process obs flow with data removed for pest run\"\n#            -> append r-script pointer on a new line within the same paragraph\n# This paragraph is the last one in the document body, and Range.InsertXML\n# cannot consume the document's final paragraph mark, so replacing it in\n# place would leave a stray empty paragraph behind. Work around that by\n# temporarily adding a throw-away paragraph after it (so our target is no\n# longer \"last\"), doing the XML replace, then deleting the throw-away\n# paragraph that is left as the new final (empty) paragraph.\n$p3 = Find-Paragraph \"Process obs flow with data removed for pest run\"\nif ($null -eq $p3) { throw \"Could not locate step 3 paragraph\" }\n$isLastParagraph = ($p3.Range.End -eq $d.Content.End)\nif ($isLastParagraph) {\n    $p3.Range.InsertParagraphAfter()\n}\n$body3 = '<w:p>' + $listPPr +\n    '<w:r><w:t xml:space=\"preserve\">Process </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>obs</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> flow with data removed for pest run</w:t></w:r>' +\n    '<w:r><w:t>, use previously developed r-script:</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:t>M:\\Models\\Bacteria\\HSPF\\bigelkHydroCal201601\\r-files</w:t></w:r>' +\n    '<w:r><w:t>\\hspf.output-proc.R</w:t></w:r>' +\n    '</w:p>'\n$p3.Range.InsertXML((New-PkgXml $body3))\nif ($isLastParagraph) {\n    $trailing = $d.Paragraphs($d.Paragraphs.Count)\n    $trailing.Range.Delete()\n}\n"}
